$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-apply the AutoFilter, pinned to A1:F35 (done before adding the new rows below
# so the filter range doesn't auto-expand to cover them)
$ws.AutoFilterMode = $false
$ws.Range("A1:F35").AutoFilter() | Out-Null

# Update the hidden _FilterDatabase defined name to match the autofilter range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$F`$35"
    }
}

# Copy the existing row's formatting pattern down into the new rows (36-39)
$src = $ws.Range("A35:F35")
for ($r = 36; $r -le 39; $r++) {
    $dst = $ws.Range("A$r`:F$r")
    $src.Copy($dst)
}

# New test rows to add (TestName / ID), rest of columns mirror existing rows
$newRows = @(
    @("Product_Summary-Credit_Card_List_[WEB]_1", "C70792"),
    @("Manage_Products-Favorite_account_[WEB]", "C70793"),
    @("Manage_Products-Favorite_account-Removal_of_the_favorite_account_[WEB]", "C70794"),
    @("Term_Deposits_Lists_[WEB]", "C70795")
)

$r = 36
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

# Move the selection to reflect where the user ended up after adding rows
$ws.Range("A41").Select() | Out-Null
